$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 90
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H19").Value = 5494883
$ws.Range("I19").Value = 11905043
$ws.Range("J19").Value = 459.64285
$ws.Range("K19").Value = 11905043
$ws.Range("L19").Value = 459.64285
$ws.Range("M19").Value = -11904868
$ws.Range("N19").Value = -809.64285

$ws.Range("H98").Value = 6645.4443
$ws.Range("I98").Value = 7829.857
$ws.Range("K98").Value = 7829.857
$ws.Range("M98").Value = -6331.857

$ws.Range("H111").Value = 556.7619
$ws.Range("I111").Value = 442.15384
$ws.Range("J111").Value = 743
$ws.Range("K111").Value = 1326.46152
$ws.Range("L111").Value = 2229
$ws.Range("M111").Value = 1740.53848
$ws.Range("N111").Value = -8363

$ws.Range("H122").Value = 6645.4443
$ws.Range("I122").Value = 7829.857
$ws.Range("K122").Value = 23489.571
$ws.Range("M122").Value = -21039.571

$ws.Range("H131").Value = 3242.9285
$ws.Range("I131").Value = 3654.889
$ws.Range("J131").Value = 2501.4
$ws.Range("K131").Value = 10964.667
$ws.Range("L131").Value = 7504.200000000001
$ws.Range("M131").Value = -5924.667000000001
$ws.Range("N131").Value = -17584.2

$ws.Range("H135").Value = 1347.625
$ws.Range("J135").Value = 919
$ws.Range("L135").Value = 8271
$ws.Range("N135").Value = -13341

$ws.Range("H137").Value = 3960.0527
$ws.Range("I137").Value = 4318.852
$ws.Range("J137").Value = 3079.3635
$ws.Range("K137").Value = 12956.556
$ws.Range("L137").Value = 9238.0905
$ws.Range("M137").Value = -10406.556
$ws.Range("N137").Value = -14338.0905

$ws.Range("H141").Value = 485910.3
$ws.Range("I141").Value = 2535.7693
$ws.Range("J141").Value = 1114297.2
$ws.Range("K141").Value = 7607.3079
$ws.Range("L141").Value = 3342891.6
$ws.Range("M141").Value = -2427.3079
$ws.Range("N141").Value = -3353251.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7577168
$ws.Range("I2").Value = 16668137
$ws.Range("K2").Value = 16668137
$ws.Range("M2").Value = -16668024

$ws.Range("H45").Value = 1184.3572
$ws.Range("I45").Value = 1083.1538
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 1083.1538
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -706.1538
$ws.Range("N45").Value = -3254

$ws.Range("H76").Value = 27371.637
$ws.Range("J76").Value = 27371.637
$ws.Range("L76").Value = 27371.637
$ws.Range("N76").Value = -28047.637

$ws.Range("H79").Value = 27371.637
$ws.Range("J79").Value = 27371.637
$ws.Range("L79").Value = 27371.637
$ws.Range("N79").Value = -29711.637

$ws.Range("H110").Value = 1795.238
$ws.Range("I110").Value = 761.1111
$ws.Range("K110").Value = 761.1111
$ws.Range("M110").Value = 1283.8889

$ws.Range("H116").Value = 7577168
$ws.Range("I116").Value = 16668137
$ws.Range("K116").Value = 16668137
$ws.Range("M116").Value = -16665843

$ws.Range("H134").Value = 32450
$ws.Range("J134").Value = 32450
$ws.Range("L134").Value = 32450
$ws.Range("N134").Value = -42590

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7577168
$ws.Range("I3").Value = 16668137
$ws.Range("K3").Value = 16668137
$ws.Range("M3").Value = -16668023

$ws.Range("H64").Value = 446.61905
$ws.Range("I64").Value = 524.7143
$ws.Range("J64").Value = 407.57144
$ws.Range("K64").Value = 524.7143
$ws.Range("L64").Value = 407.57144
$ws.Range("M64").Value = -299.7143
$ws.Range("N64").Value = -857.5714399999999

$ws.Range("H67").Value = 446.61905
$ws.Range("I67").Value = 524.7143
$ws.Range("J67").Value = 407.57144
$ws.Range("K67").Value = 524.7143
$ws.Range("L67").Value = 407.57144
$ws.Range("M67").Value = 255.2857
$ws.Range("N67").Value = -1967.57144

$ws.Range("H86").Value = 75522.71000000001
$ws.Range("J86").Value = 102755.8
$ws.Range("L86").Value = 102755.8
$ws.Range("N86").Value = -105001.8

$ws.Range("H89").Value = 75522.71000000001
$ws.Range("J89").Value = 102755.8
$ws.Range("L89").Value = 513779
$ws.Range("N89").Value = -525011

$ws.Range("H107").Value = 990.36365
$ws.Range("I107").Value = 939.625
$ws.Range("K107").Value = 939.625
$ws.Range("M107").Value = 980.375

$ws.Range("H134").Value = 2285.3103
$ws.Range("I134").Value = 2433.8367
$ws.Range("K134").Value = 7301.5101
$ws.Range("M134").Value = -4766.5101

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2348.4565
$ws.Range("I31").Value = 1670.6757
$ws.Range("J31").Value = 5134.8887
$ws.Range("K31").Value = 1670.6757
$ws.Range("L31").Value = 5134.8887
$ws.Range("M31").Value = -1375.6757
$ws.Range("N31").Value = -5724.8887

$ws.Range("H34").Value = 2348.4565
$ws.Range("I34").Value = 1670.6757
$ws.Range("J34").Value = 5134.8887
$ws.Range("K34").Value = 1670.6757
$ws.Range("L34").Value = 5134.8887
$ws.Range("M34").Value = -1468.6757
$ws.Range("N34").Value = -5538.8887

$ws.Range("H58").Value = 11366576
$ws.Range("I58").Value = 1984.1428
$ws.Range("J58").Value = 31254612
$ws.Range("K58").Value = 1984.1428
$ws.Range("L58").Value = 31254612
$ws.Range("M58").Value = -1781.1428
$ws.Range("N58").Value = -31255018

$ws.Range("H99").Value = 1897.5714
$ws.Range("I99").Value = 1732.5
$ws.Range("J99").Value = 2117.6667
$ws.Range("K99").Value = 1732.5
$ws.Range("L99").Value = 2117.6667
$ws.Range("M99").Value = -234.5
$ws.Range("N99").Value = -5113.6667

$ws.Range("H107").Value = 603.1613
$ws.Range("I107").Value = 689.4286
$ws.Range("K107").Value = 689.4286
$ws.Range("M107").Value = 1230.5714

$ws.Range("H126").Value = 1897.5714
$ws.Range("I126").Value = 1732.5
$ws.Range("J126").Value = 2117.6667
$ws.Range("K126").Value = 5197.5
$ws.Range("L126").Value = 6353.000100000001
$ws.Range("M126").Value = -2727.5
$ws.Range("N126").Value = -11293.0001

$ws.Range("H132").Value = 2409.5454
$ws.Range("I132").Value = 3591.3
$ws.Range("J132").Value = 1424.75
$ws.Range("K132").Value = 10773.9
$ws.Range("L132").Value = 4274.25
$ws.Range("M132").Value = -8243.900000000001
$ws.Range("N132").Value = -9334.25

$ws.Range("H134").Value = 3006.814
$ws.Range("I134").Value = 3267.258
$ws.Range("J134").Value = 2334
$ws.Range("K134").Value = 9801.773999999999
$ws.Range("L134").Value = 7002
$ws.Range("M134").Value = -7266.773999999999
$ws.Range("N134").Value = -12072

$ws.Range("H136").Value = 11366576
$ws.Range("I136").Value = 1984.1428
$ws.Range("J136").Value = 31254612
$ws.Range("K136").Value = 5952.428400000001
$ws.Range("L136").Value = 93763836
$ws.Range("M136").Value = -3402.428400000001
$ws.Range("N136").Value = -93768936

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 953.1429000000001
$ws.Range("J47").Value = 2950
$ws.Range("L47").Value = 8850
$ws.Range("N47").Value = -9712

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 288721.72
$ws.Range("I126").Value = 2339.2
$ws.Range("J126").Value = 403274.72
$ws.Range("K126").Value = 7017.599999999999
$ws.Range("L126").Value = 1209824.16
$ws.Range("M126").Value = -4547.599999999999
$ws.Range("N126").Value = -1214764.16

$ws.Range("H132").Value = 3393.5952
$ws.Range("J132").Value = 2991.7827
$ws.Range("L132").Value = 8975.348100000001
$ws.Range("N132").Value = -14035.3481

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2131.9167
$ws.Range("I61").Value = 564.7778
$ws.Range("K61").Value = 564.7778
$ws.Range("M61").Value = -362.7778

$ws.Range("H113").Value = 2131.9167
$ws.Range("I113").Value = 564.7778
$ws.Range("K113").Value = 564.7778
$ws.Range("M113").Value = 1605.2222

$ws.Range("H136").Value = 3862.6875
$ws.Range("I136").Value = 3982.0908
$ws.Range("J136").Value = 3600
$ws.Range("K136").Value = 11946.2724
$ws.Range("L136").Value = 10800
$ws.Range("M136").Value = -9396.2724
$ws.Range("N136").Value = -15900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 771.4
$ws.Range("I113").Value = 405.3684
$ws.Range("K113").Value = 1216.1052
$ws.Range("M113").Value = 953.8948

$ws.Range("H122").Value = 2459.4
$ws.Range("I122").Value = 2066
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 6198
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -3748
$ws.Range("N122").Value = -22900

$ws.Range("H132").Value = 5055.4614
$ws.Range("I132").Value = 2176.1292
$ws.Range("K132").Value = 6528.3876
$ws.Range("M132").Value = -3998.3876

$ws.Range("H136").Value = 3173.2222
$ws.Range("I136").Value = 3210.0833
$ws.Range("K136").Value = 9630.249899999999
$ws.Range("M136").Value = -7080.249899999999
